$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2:N5").Value = 85.8724807945396
